$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Sheet1"
$ws.Range("F21").Select()
